$d = $word.ActiveDocument

# 1) Merge the three hyperlink runs ("https://www" + "." + "deapoio.com.br/cursos/chama.php")
#    into a single run with the full URL text.
$h = $d.Hyperlinks(1)
$h.TextToDisplay = "https://www.deapoio.com.br/cursos/chama.php"

# 2) Remove the entire "Processo II - Só Procedural" section (flowchart + its
#    accompanying pseudo-code/parametrization text) that used to follow the
#    "Não completei essa funcionalidade..." paragraph, all the way to the end
#    of the document body (right before the final sectPr).
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $txt = $d.Paragraphs($i).Range.Text
    if ($txt -like "*completei essa funcionalidade*") {
        $cutStart = $d.Paragraphs($i + 1).Range.Start
        break
    }
}
$cutEnd = $d.Content.End
$r = $d.Range($cutStart, $cutEnd)
$r.Delete()
